# "Generate Report for Archive"
# Localization status moved from "Ready for handoff" to "In Translation"
# for the two source files tracked in this workbook. Update every cell
# that carries that status (Overview summary columns + each per-locale
# "Status" column), then tighten the now-narrower Status columns to match.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (col E) / de-de (col F) status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Status text shrank, so the columns re-fit narrower.
$wsOverview.Columns.Item(5).ColumnWidth = 13
$wsOverview.Columns.Item(6).ColumnWidth = 13

# --- zh-cn sheet: Status column (col C) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus
$wsZh.Columns.Item(3).ColumnWidth = 13

# --- de-de sheet: Status column (col C) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus
$wsDe.Columns.Item(3).ColumnWidth = 13
